# The "建物" (building) sheet's property_category column (I) was
# mislabeled "land" for every row -- correct it to "building".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")
$ws.Range("I2:I8").Value = "building"
